# Update the speaker notes on the slide with sldId="320" (13th slide in the
# presentation's slide order) to add a bolded "Project Challenges" heading,
# split the existing narrative into separate paragraphs, tweak some of the
# wording, and append a closing "grades" paragraph.

$p = $ppt.ActivePresentation

# sldId="320" is the 13th slide in the presentation's slide order (p:sldIdLst).
$s = $p.Slides.Item(13)

$notesShape = $s.NotesPage.Shapes.Item(2)
$tr = $notesShape.TextFrame.TextRange

$line1 = "Project Challenges"
$line2 = "I was required to complete multiple concurrent tasks, some with unfamiliar tools, within a strict time limit."
$line3 = "Deliverables were; detailing accurate results in three technical reports, plus producing and delivering a technical presentation."
$line4 = "I needed to:  Select the appropriate tools at each stage, Execute with precision, Document my progress, Produce accurate and clear reports, Implement targeted monitoring with alerts, and Recommend solutions for the vulnerabilities found.  "
$line5 = "My project submission earned an A- and my overall academic average was A+"

$tr.Text = $line1 + "`n" + $line2 + "`n" + $line3 + "`n" + $line4 + "`n" + $line5

# Make the new heading paragraph bold.
$tr.Paragraphs(1,1).Font.Bold = $true

# Match the smaller font size used by the "I needed to..." / grades paragraphs.
$tr.Paragraphs(4,1).Font.Size = 12
$tr.Paragraphs(5,1).Font.Size = 12
